$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.769.16"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "3.688.64"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.88"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("E6").Value = "  +15.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "669.63"
$ws.Range("E7").Value = "  +2.47%  "
$ws.Range("E8").Value = "  +4.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.10"
$ws.Range("E9").Value = "  +4.86%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "3.687.69"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.55"
$ws.Range("E12").Value = "  +5.11%  "
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.62"
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("D15").Value = "4.372.29"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("D17").Value = "96.488.18"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.04"
$ws.Range("E18").Value = "  +16.63%  "
$ws.Range("D19").Value = "3.688.20"
$ws.Range("E19").Value = "  +3.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.82"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.46"
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.534"
$ws.Range("E22").Value = "  +3.73%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "520.97"
$ws.Range("E23").Value = "  +3.27%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.50"
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("E25").Value = "  +4.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.02"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.90"
$ws.Range("E27").Value = "  +7.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.08"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  +7.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.06"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.20"
$ws.Range("E31").Value = "  +7.00%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +10.81%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.92"
$ws.Range("E35").Value = "  +4.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.588"
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "623.10"
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.81"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.76"
$ws.Range("E40").Value = "  +28.81%  "
$ws.Range("E41").Value = "  +6.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.960"
$ws.Range("E42").Value = "  +6.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.96"
$ws.Range("E43").Value = "  +7.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.21"
$ws.Range("E45").Value = "  +8.53%  "
$ws.Range("E46").Value = "  +7.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.432"
$ws.Range("E47").Value = "  +25.05%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("E50").Value = "  +4.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.71"
$ws.Range("E51").Value = "  +3.73%  "
